# Updated symbol list on Mon Dec 19 14:40:25 UTC 2022 with GitHub Actions
# Refreshes the crypto price/volume table: updates prices in column D,
# shifts the WazirX..One coin block down one row (inserting "One" at the
# top of that block) and refreshes a few "Best/Worst in 24h" suffixes in
# column E.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Address, $Text)
    $cell = $ws.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
}

Set-TextValue 'D2' '249.33'
Set-TextValue 'D3' '21.85'
Set-TextValue 'D4' '5.437'
Set-TextValue 'D5' '0.05698'
Set-TextValue 'D6' '3.379'
Set-TextValue 'D8' '1.025'
$ws.Range('B9').Value = 'One'
$ws.Range('C9').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue 'D9' '0.0005865'
$ws.Range('E9').Value = '8OneONE'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue 'D10' '0.1453'
$ws.Range('E10').Value = '9WazirXWRX'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue 'D11' '0.07608'
$ws.Range('E11').Value = '10MandalaExchangeTokenMDX'
$ws.Range('B12').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C12').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue 'D12' '0.03167'
$ws.Range('E12').Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue 'D13' '0.03029'
$ws.Range('E13').Value = '12BitrueCoinBTR'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue 'D14' '0.09262'
$ws.Range('E14').Value = '13BitMartTokenBMX'
$ws.Range('B15').Value = 'MCDex'
$ws.Range('C15').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue 'D15' '3.532'
$ws.Range('E15').Value = '14MCDexMCB'
$ws.Range('B16').Value = 'BitForexToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue 'D16' '0.001663'
$ws.Range('E16').Value = '15BitForexTokenBF'
$ws.Range('B17').Value = 'CoinExToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue 'D17' '0.04719'
$ws.Range('E17').Value = '16CoinExTokenCET'
Set-TextValue 'D18' '0.006352'
$ws.Range('E19').Value = '18HotbitTokenHTB'
Set-TextValue 'D21' '0.0001502'
Set-TextValue 'D22' '0.0003104'
Set-TextValue 'D23' '3.768'
Set-TextValue 'D24' '6.421'
Set-TextValue 'D25' '2.158'
Set-TextValue 'D26' '0.3295'
Set-TextValue 'D27' '0.1336'
Set-TextValue 'D40' '0.04066'
Set-TextValue 'D41' '0.006982'
$ws.Range('E41').Value = '40KickTokenKICKBestin24h'
Set-TextValue 'D42' '0.003505'
Set-TextValue 'D43' '0.1040'
Set-TextValue 'D45' '0.00005911'
Set-TextValue 'D47' '0.0005504'
$ws.Range('E47').Value = '46ACDXExchangeACXTWorstin24h'
Set-TextValue 'D48' '0.6833'
Set-TextValue 'D49' '0.008038'
Set-TextValue 'D50' '0.00002103'
